# Adds two new "Title and Content" slides (Ongoing processes / In an
# emergency) after the existing slide, per the authored diff.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2: "Ongoing processes"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)

$s2Title = $s2.Shapes.Item(1)
$s2Title.Left = 56.742916107177734
$s2Title.Top = 0.0
$s2Title.Width = 828.0
$s2Title.Height = 80.22858428955078
$s2Title.TextFrame.TextRange.Text = 'Ongoing processes'

$s2Body = $s2.Shapes.Item(2)
$s2Body.Left = 66.0
$s2Body.Top = 80.22858428955078
$s2Body.Width = 828.0
$s2Body.Height = 406.146484375

$s2tr = $s2Body.TextFrame.TextRange
$s2tr.Text = '‘nudge’ for preparedness, eg.'
[void]$s2tr.InsertAfter("`nGutter cleaning")
[void]$s2tr.InsertAfter("`nPump/sprinkler testing")
[void]$s2tr.InsertAfter("`nReconnaissance of potential escape routes")
[void]$s2tr.InsertAfter("`nGeocaching as an incentive (gamification) mean")
[void]$s2tr.InsertAfter("`nSocial media integration to help build relationships and plans with neighbours, eg. Awareness of others’")
[void]$s2tr.InsertAfter("`nHabits")
[void]$s2tr.InsertAfter("`nStrategies")
[void]$s2tr.InsertAfter("`nLimitations (eg. physical impairment, dependents: children, elderly, pets, etc). In an emergency situation connections with neighbours may be used to help nudge safe behaviour (eg. " + [char]0x201C + "90% of your neighbours have already evacuated" + [char]0x201D + ")")

$s2Levels = @{
    2 = 2
    3 = 2
    4 = 2
    5 = 3
    7 = 2
    8 = 2
    9 = 2
}
foreach ($idx in $s2Levels.Keys) {
    $s2tr.Paragraphs($idx, 1).IndentLevel = $s2Levels[$idx]
}

# ---------------------------------------------------------------------
# Slide 3: "In an emergency"
# ---------------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)

$s3Title = $s3.Shapes.Item(1)
$s3Title.Left = 56.742916107177734
$s3Title.Top = 0.0
$s3Title.Width = 828.0
$s3Title.Height = 80.22858428955078
$s3Title.TextFrame.TextRange.Text = 'In an emergency'

$s3Body = $s3.Shapes.Item(2)
$s3Body.Left = 66.0
$s3Body.Top = 80.22858428955078
$s3Body.Width = 828.0
$s3Body.Height = 406.146484375

$s3tr = $s3Body.TextFrame.TextRange
$s3tr.Text = 'Update escape routes advice as situation evolves'
[void]$s3tr.InsertAfter("`nProvide advice on likely arrival time of fire and suggest action deadlines")

Write-Output ("Slides now: " + $p.Slides.Count)
